# "Injection table - precalculated" — update the raw injection-multiplier
# grid (rows 27:36) so that the values forming the upper "ramp" become 1,
# matching the rest of the already-1 cells in each of those rows. The
# downstream ROUND(...) formulas in rows 41:55 are driven off this table
# and recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = 1

$ws.Range("A28:B28").Value = 1

$ws.Range("A29:C29").Value = 1

$ws.Range("A30:C30").Value = 1

$ws.Range("A31:E31").Value = 1
$ws.Range("A32:E32").Value = 1
$ws.Range("A33:E33").Value = 1
$ws.Range("A34:E34").Value = 1
$ws.Range("A35:E35").Value = 1
$ws.Range("A36:E36").Value = 1

# Make sure the dependent formulas (ROUND-based table in rows 41:55) are
# fully recalculated with the new inputs before the workbook is saved.
$excel.CalculateFullRebuild()

# Restore the view/selection state the sheet was left in: scrolled so
# row 4 is at the top, with H28 selected (previously topLeftCell A22 /
# A40:P55 selected).
$win = $wb.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("H28").Select()
